$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.227.86"
$ws.Range("E2").Value = "  +0.58%  "

$ws.Range("D3").Value = "'1.895.98"
$ws.Range("E3").Value = "  +0.10%  "

$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").Value = "'307.33"
$ws.Range("E5").Value = "  +0.08%  "

$ws.Range("E6").Value = "  +0.09%  "

$ws.Range("D7").Value = "'0.5189"
$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").Value = "'0.3761"
$ws.Range("E8").Value = "  -0.22%  "

$ws.Range("D9").Value = "'0.07282"
$ws.Range("E9").Value = "  +0.80%  "

$ws.Range("D10").Value = "'21.19"
$ws.Range("E10").Value = "  +0.24%  "

$ws.Range("D11").Value = "'0.8999"
$ws.Range("E11").Value = "  +0.89%  "

$ws.Range("D12").Value = "'0.08156"
$ws.Range("E12").Value = "  +6.42%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "'1.912.73"
$ws.Range("E13").Value = "  +1.02%  "

$ws.Range("B14").Value = "Litecoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D14").Value = "'96.38"
$ws.Range("E14").Value = "  +2.32%  "

$ws.Range("D15").Value = "'5.276"
$ws.Range("E15").Value = "  +0.97%  "

$ws.Range("D16").Value = "'1.002"
$ws.Range("E16").Value = "  +0.08%  "

$ws.Range("D17").Value = "'0.000008620"
$ws.Range("E17").Value = "  +1.31%  "

$ws.Range("E18").Value = "  +0.22%  "

$ws.Range("E19").Value = "  +0.13%  "

$ws.Range("D20").Value = "'27.258.71"
$ws.Range("E20").Value = "  +0.49%  "

$ws.Range("D21").Value = "'5.083"
$ws.Range("E21").Value = "  +0.41%  "

$ws.Range("D22").Value = "'10.68"
$ws.Range("E22").Value = "  +0.99%  "

$ws.Range("D23").Value = "'6.393"
$ws.Range("E23").Value = "  -0.29%  "

$ws.Range("D24").Value = "'2.298"
$ws.Range("E24").Value = "  +0.28%  "

$ws.Range("E25").Value = "  +0.59%  "

$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").Value = "'18.21"
$ws.Range("E26").Value = "  +0.84%  "

$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D27").Value = "'1.744"
$ws.Range("E27").Value = "  +0.59%  "

$ws.Range("D28").Value = "'115.13"
$ws.Range("E28").Value = "  +0.59%  "

$ws.Range("D29").Value = "'4.970"
$ws.Range("E29").Value = "  +0.12%  "

$ws.Range("D30").Value = "'4.835"
$ws.Range("E30").Value = "  +1.03%  "

$ws.Range("D31").Value = "'0.09231"
$ws.Range("E31").Value = "  +0.38%  "

$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "'0.7952"
$ws.Range("E32").Value = "  +2.47%  "

$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.05033"
$ws.Range("E33").Value = "  -0.27%  "

$ws.Range("D34").Value = "'1.215"
$ws.Range("E34").Value = "  -1.65%  "

$ws.Range("D35").Value = "'3.451"
$ws.Range("E35").Value = "  +4.75%  "

$ws.Range("D36").Value = "'2.952"
$ws.Range("E36").Value = "  -0.78%  "

$ws.Range("D37").Value = "'2.600"
$ws.Range("E37").Value = "  +0.69%  "

$ws.Range("D38").Value = "'0.5667"
$ws.Range("E38").Value = "  +0.94%  "

$ws.Range("D39").Value = "'0.01986"
$ws.Range("E39").Value = "  -0.11%  "

$ws.Range("E40").Value = "  +0.09%  "

$ws.Range("D41").Value = "'8.953"
$ws.Range("E41").Value = "  -0.35%  "

$ws.Range("D42").Value = "'6.557"
$ws.Range("E42").Value = "  -1.04%  "

$ws.Range("D43").Value = "'115.29"
$ws.Range("E43").Value = "  -2.86%  "

$ws.Range("D44").Value = "'0.1514"
$ws.Range("E44").Value = "  -0.17%  "

$ws.Range("D45").Value = "'0.4854"
$ws.Range("E45").Value = "  +0.73%  "

$ws.Range("E46").Value = "  +0.08%  "

$ws.Range("E47").Value = "  -1.42%  "

$ws.Range("D48").Value = "'1.620"
$ws.Range("E48").Value = "  +1.68%  "

$ws.Range("D49").Value = "'38.18"
$ws.Range("E49").Value = "  +1.82%  "

$ws.Range("D50").Value = "'63.38"
$ws.Range("E50").Value = "  -0.86%  "

$ws.Range("D51").Value = "'0.05941"
$ws.Range("E51").Value = "  +0.29%  "
